# Re-generated project files: update PID4Cat Excel data-entry model.
$wb = $excel.ActiveWorkbook

# --- Sheet "PID4CatRecord": drop the `record_version` column and rename a
#     few headers (D:J -> D:I after the shift). ---
$wsRecord = $wb.Worksheets.Item("PID4CatRecord")

# Remove column D ("record_version"); everything to the right shifts left,
# so the sheet goes from A1:J1 to A1:I1 automatically.
$wsRecord.Columns.Item(4).Delete(-4159)   # xlShiftToLeft

# After the shift: D=pid_schema_version, E=dc_rights, F=curation_contact,
# G=resource_info, H=related_identifiers, I=change_log.
# Rename the two fields that also changed name in this revision.
$wsRecord.Range("E1").Value = "license"
$wsRecord.Range("F1").Value = "curation_contact_email"

# --- Sheet "ResourceInfo": extend the resource_category dropdown list. ---
$wsResourceInfo = $wb.Worksheets.Item("ResourceInfo")
$wsResourceInfo.Range("C2:C1048576").Validation.Formula1 = '"COLLECTION,SAMPLE,MATERIAL,DEVICE,DATA_OBJECT,DATA_SERVICE"'

# --- Sheet "LogRecord": swap RIGHTS for LICENSE in the changed_field list. ---
$wsLogRecord = $wb.Worksheets.Item("LogRecord")
$wsLogRecord.Range("C2:C1048576").Validation.Formula1 = '"STATUS,RESOURCE_INFO,RELATED_IDS,CONTACT,LICENSE"'

# --- Sheet "Agent": rename contact columns. ---
$wsAgent = $wb.Worksheets.Item("Agent")
$wsAgent.Range("B1").Value = "email"
$wsAgent.Range("C1").Value = "orcid"
